$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the duplicate "Adrenal adenoma / Clip 2 B-mode" row (row 14) -
# it duplicated row 13's Adrenal adenoma entry with a placeholder clip.
$ws.Rows.Item(14).Delete()

# Remove the leftover test row at the bottom of the table
# (originally row 18, now row 17 after the deletion above).
$ws.Rows.Item(17).Delete()

# Re-apply the existing sort (Organ column, no header) over the new,
# smaller data range so the sheet's remembered sort state matches it.
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("A2:A16"))
$ws.Sort.SetRange($ws.Range("A2:C16"))
$ws.Sort.Header = 0
$ws.Sort.Apply()

# Leave the selection where the author left it when saving.
$ws.Range("B11").Select()
